$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rich-text edits: change only the digits inside existing runs,
# preserving the surrounding run text of the shared strings.
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "52"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "12/23/2024"
$c9.Characters(48, 10).Text = "12/29/2024"

# --- Cells that must hold the literal text "0" (shared string), not the
# number 0. Excel auto-converts numeric-looking text back to a number
# unless the cell is explicitly formatted as Text first.
foreach ($ref in @("D15", "C22", "D22", "D27")) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value2 = "0"
}

# --- Remaining numeric / text value updates in the crime-stats table
$ws.Range("L14").Value2 = -41.666666666666
$ws.Range("C15").Value2 = 1
$ws.Range("E15").Value2 = "***.*"
$ws.Range("G15").Value2 = 4
$ws.Range("H15").Value2 = -75
$ws.Range("I15").Value2 = 41
$ws.Range("K15").Value2 = 36.666666666666
$ws.Range("L15").Value2 = 20.588235294117
$ws.Range("M15").Value2 = 95.238095238095
$ws.Range("N15").Value2 = -44.594594594594
$ws.Range("C16").Value2 = 9
$ws.Range("D16").Value2 = 15
$ws.Range("E16").Value2 = -40
$ws.Range("F16").Value2 = 30
$ws.Range("G16").Value2 = 41
$ws.Range("H16").Value2 = -26.829268292682
$ws.Range("I16").Value2 = 437
$ws.Range("J16").Value2 = 387
$ws.Range("K16").Value2 = 12.919896640826
$ws.Range("L16").Value2 = 13.212435233160
$ws.Range("M16").Value2 = -9.147609147609
$ws.Range("N16").Value2 = -74.607786170830
$ws.Range("C17").Value2 = 13
$ws.Range("D17").Value2 = 8
$ws.Range("E17").Value2 = 62.5
$ws.Range("F17").Value2 = 54
$ws.Range("G17").Value2 = 39
$ws.Range("H17").Value2 = 38.461538461538
$ws.Range("I17").Value2 = 780
$ws.Range("J17").Value2 = 653
$ws.Range("K17").Value2 = 19.448698315467
$ws.Range("L17").Value2 = 37.082601054481
$ws.Range("M17").Value2 = 128.739002932551
$ws.Range("N17").Value2 = -7.582938388625
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 25
$ws.Range("G18").Value2 = 14
$ws.Range("H18").Value2 = 78.571428571428
$ws.Range("I18").Value2 = 234
$ws.Range("J18").Value2 = 164
$ws.Range("K18").Value2 = 42.682926829268
$ws.Range("L18").Value2 = 11.961722488038
$ws.Range("M18").Value2 = -9.652509652509
$ws.Range("N18").Value2 = -79.792746113989
$ws.Range("C19").Value2 = 17
$ws.Range("D19").Value2 = 6
$ws.Range("E19").Value2 = 183.333333333333
$ws.Range("F19").Value2 = 39
$ws.Range("G19").Value2 = 38
$ws.Range("H19").Value2 = 2.631578947368
$ws.Range("I19").Value2 = 592
$ws.Range("J19").Value2 = 582
$ws.Range("K19").Value2 = 1.718213058419
$ws.Range("L19").Value2 = -1.168614357262
$ws.Range("M19").Value2 = 39.952718676122
$ws.Range("N19").Value2 = -50.871369294605
$ws.Range("C20").Value2 = 4
$ws.Range("E20").Value2 = 100
$ws.Range("F20").Value2 = 18
$ws.Range("G20").Value2 = 12
$ws.Range("H20").Value2 = 50
$ws.Range("I20").Value2 = 201
$ws.Range("J20").Value2 = 226
$ws.Range("K20").Value2 = -11.061946902654
$ws.Range("L20").Value2 = -4.285714285714
$ws.Range("M20").Value2 = 18.235294117647
$ws.Range("N20").Value2 = -87.286527514231
$ws.Range("C21").Value2 = 48
$ws.Range("D21").Value2 = 35
$ws.Range("E21").Value2 = 37.142857142857
$ws.Range("F21").Value2 = 167
$ws.Range("G21").Value2 = 148
$ws.Range("H21").Value2 = 12.837837837837
$ws.Range("I21").Value2 = 2292
$ws.Range("J21").Value2 = 2044
$ws.Range("K21").Value2 = 12.133072407045
$ws.Range("L21").Value2 = 13.521545319465
$ws.Range("M21").Value2 = 34.113516676419
$ws.Range("N21").Value2 = -65.330509756466
$ws.Range("E22").Value2 = "***.*"
$ws.Range("G22").Value2 = 1
$ws.Range("H22").Value2 = 0
$ws.Range("L22").Value2 = 24
$ws.Range("C23").Value2 = 2
$ws.Range("E23").Value2 = 0
$ws.Range("F23").Value2 = 4
$ws.Range("G23").Value2 = 4
$ws.Range("I23").Value2 = 44
$ws.Range("J23").Value2 = 45
$ws.Range("K23").Value2 = -2.222222222222
$ws.Range("L23").Value2 = -15.384615384615
$ws.Range("M23").Value2 = 18.918918918918
$ws.Range("C24").Value2 = 32
$ws.Range("D24").Value2 = 36
$ws.Range("E24").Value2 = -11.111111111111
$ws.Range("F24").Value2 = 168
$ws.Range("G24").Value2 = 156
$ws.Range("H24").Value2 = 7.692307692307
$ws.Range("I24").Value2 = 2110
$ws.Range("J24").Value2 = 1616
$ws.Range("K24").Value2 = 30.569306930693
$ws.Range("L24").Value2 = 43.929058663028
$ws.Range("M24").Value2 = 88.561215370866
$ws.Range("C25").Value2 = 17
$ws.Range("D25").Value2 = 12
$ws.Range("E25").Value2 = 41.666666666666
$ws.Range("F25").Value2 = 94
$ws.Range("H25").Value2 = 38.235294117647
$ws.Range("I25").Value2 = 1310
$ws.Range("J25").Value2 = 685
$ws.Range("K25").Value2 = 91.240875912408
$ws.Range("L25").Value2 = 94.650817236255
$ws.Range("C26").Value2 = 13
$ws.Range("D26").Value2 = 19
$ws.Range("E26").Value2 = -31.578947368421
$ws.Range("F26").Value2 = 63
$ws.Range("G26").Value2 = 57
$ws.Range("H26").Value2 = 10.526315789473
$ws.Range("I26").Value2 = 1068
$ws.Range("J26").Value2 = 964
$ws.Range("K26").Value2 = 10.788381742738
$ws.Range("L26").Value2 = 36.398467432950
$ws.Range("M26").Value2 = 38.701298701298
$ws.Range("C27").Value2 = 1
$ws.Range("E27").Value2 = "***.*"
$ws.Range("F27").Value2 = 3
$ws.Range("G27").Value2 = 6
$ws.Range("H27").Value2 = -50
$ws.Range("I27").Value2 = 60
$ws.Range("K27").Value2 = 39.534883720930
$ws.Range("L27").Value2 = 22.448979591836
$ws.Range("C28").Value2 = 1
$ws.Range("D28").Value2 = 3
$ws.Range("E28").Value2 = -66.666666666666
$ws.Range("F28").Value2 = 6
$ws.Range("H28").Value2 = -14.285714285714
$ws.Range("I28").Value2 = 118
$ws.Range("J28").Value2 = 91
$ws.Range("K28").Value2 = 29.670329670329
$ws.Range("L28").Value2 = 34.090909090909
$ws.Range("L29").Value2 = -24.242424242424
$ws.Range("N29").Value2 = -84.567901234567
$ws.Range("L30").Value2 = -28.571428571428
$ws.Range("N30").Value2 = -86.301369863013
